$d = $word.ActiveDocument

# =====================================================================
# Change 1: "For more details refer to the ArchitectureGuide.docx, that
#           is part of the transfer package. The following components
#           were implemented." paragraph.
#
#   - Split "For more details refer to the " into three runs:
#       "For more details" | " about the system components " | "refer to the "
#   - Move the "_GoBack" bookmark to sit right after "ArchitectureGuide.docx"
#     (Word keeps only one "_GoBack" bookmark per document, so adding it
#     here automatically removes it from its old location near "Fixes").
#   - Replace the trailing run's text.
# =====================================================================

$lead = $d.Content
$lead.Find.Execute("For more details refer to the ", $true, $false, $false, $false, $false,
                    $true, 1, $false, "", 0) | Out-Null
$leadStart = $lead.Start
$leadEnd = $lead.End
$b1 = $leadStart + 16   # boundary between "For more details" and " "
$b2 = $leadStart + 17   # boundary between " " and "refer to the "

$arch = $d.Content
$arch.Find.Execute("ArchitectureGuide.docx", $true, $false, $false, $false, $false,
                    $true, 1, $false, "", 0) | Out-Null
$archEnd = $arch.End

# Process right-to-left so earlier offsets stay valid (ranges are plain
# fixed offsets here and do not auto-shift when other text changes length).

# 1) Move the permanent _GoBack bookmark to sit right after "ArchitectureGuide.docx".
#    Doing this first means it also acts as a hard separator so the tail-text edit
#    below can never coalesce back into the "ArchitectureGuide.docx" run.
$d.Bookmarks.Add("_GoBack", $d.Range($archEnd, $archEnd))

# 2) Replace the tail text (rightmost text edit). Set a throwaway value first so the
#    final assignment is a genuine change and xml:space/rPr get recomputed cleanly
#    instead of being left over from the original (longer) run.
$tail = $d.Content
$tail.Find.Execute(", that is part of the transfer package. The following components were implemented.",
                    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tailStart = $tail.Start
$tail.Text = "zzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"
$tail2 = $d.Range($tailStart, $tailStart + 44)
$tail2.Text = ". The following components were implemented."

# 3) Force a run split at b2 (between " " and "refer to the "), then set that piece's text.
$d.Bookmarks.Add("ZZtemp2", $d.Range($b2, $b2))
$p3 = $d.Range($b2, $leadEnd)
$p3.Text = "refer to the "

# 4) Force a run split at b1 (between "For more details" and " "), then set the middle
#    piece's text to the new inserted phrase.
$d.Bookmarks.Add("ZZtemp1", $d.Range($b1, $b1))
$p2 = $d.Range($b1, $b2)
$p2.Text = " about the system components "

# 5) Explicitly (re)set the first piece's text too, so it is not left as an untouched
#    "stale" run fragment inheriting the original run's xml:space="preserve". Go via a
#    throwaway value first so the final assignment is a genuine change.
$p1 = $d.Range($leadStart, $b1)
$p1.Text = "zzzzzzzzzzzzzzzz"
$p1 = $d.Range($leadStart, $leadStart + 16)
$p1.Text = "For more details"

# Clean up temporary split-marker bookmarks (keep _GoBack).
$d.Bookmarks("ZZtemp1").Delete()
$d.Bookmarks("ZZtemp2").Delete()

# =====================================================================
# Change 2 & 3: merge the " " and "loading module" runs that follow
#               "Registrations" / "Surveys" into a single run " loading module".
# =====================================================================

foreach ($word0 in @("Registrations", "Surveys")) {
    $anchor = $d.Content
    $anchor.Find.Execute($word0, $true, $false, $false, $false, $false,
                          $true, 1, $false, "", 0) | Out-Null
    $anchorEnd = $anchor.End

    $m = $d.Range($anchorEnd, $d.Content.End)
    $m.Find.Execute(" loading module", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0) | Out-Null
    $mStart = $m.Start

    # Protect the preceding run (e.g. "Registrations") from being swept into the
    # coalesce pass triggered by the text edits below.
    $d.Bookmarks.Add("ZZsep", $d.Range($anchorEnd, $anchorEnd))

    $m.Text = "zzzzzzzzzzzzzzz"
    $m2 = $d.Range($mStart, $mStart + 15)
    $m2.Text = " loading module"

    $d.Bookmarks("ZZsep").Delete()
}

Write-Output "done"
